$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values for rows 5-10 (empadronador / total_registros) according to the diff:
# Row5: URBINA ANTICONA ALEX BRUNO, 97
# Row6: SANCHEZ SALDAÑA FRANK REGINALDO, 97
# Row7: BLANCO LOZANO ANDREA MILAGROS, 85
# Row8: BURGA MEDINA SHIRLEY ROCIO, 85
# Row9: DE LA CRUZ CARDENAS RUTH LUCERO, 84
# Row10: LLANOS HUACCHA BRITSY, 81

$ws.Range("A5").Value = "URBINA ANTICONA ALEX BRUNO"
$ws.Range("B5").Value = 97

$ws.Range("A6").Value = "SANCHEZ SALDAÑA FRANK REGINALDO"
$ws.Range("B6").Value = 97

$ws.Range("A7").Value = "BLANCO LOZANO ANDREA MILAGROS"
$ws.Range("B7").Value = 85

$ws.Range("A8").Value = "BURGA MEDINA SHIRLEY ROCIO"
$ws.Range("B8").Value = 85

$ws.Range("A9").Value = "DE LA CRUZ CARDENAS RUTH LUCERO"
$ws.Range("B9").Value = 84

$ws.Range("A10").Value = "LLANOS HUACCHA BRITSY"
$ws.Range("B10").Value = 81
